$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New log entry for 10/21 (row 29) ---
# Date
$ws.Range("B29").Value = 45218
# Start time 10:10 AM / End time 11:01 AM
$ws.Range("C29").Value = 0.4236111111111111
$ws.Range("D29").Value = 0.45902777777777781
$ws.Range("D29").NumberFormat = $ws.Range("D28").NumberFormat

# "What I accomplished" / "What do I think I should do next session" for the new entry.
# Set the brand-new "Got backend auth started" text first so new shared strings are
# appended to the shared-string table in the same order the author typed them.
$ws.Range("G29").Value = "Got backend auth started"

# Update the prior entry's (row 28) accomplishment text - it had said "Got Oauth
# started" but that became "Got Oauth working on the frontend" (the "next session"
# note about wiring up the database/credential side stayed the same in H28).
$ws.Range("G28").Value = "Got Oauth working on the frontend"

$ws.Range("H29").Value = "Next is finishing backend auth"

# Leave the cursor where the author finished up - on H30, scrolled down a bit.
$ws.Activate()
$ws.Range("H30").Select()
$excel.ActiveWindow.ScrollRow = 26
